$d = $word.ActiveDocument
$d.Content.Find.Execute("Summary 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Summary 0", 2)
